# Aggiornamento fino a 20/09/2021
# Append new daily rows (375-385) to the existing data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 374

# New data: date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti
$newData = @(
    @(44449, 0, 0, 0),
    @(44450, 0, 0, 0),
    @(44451, 0, 0, 0),
    @(44452, 1, 1, 37.46721618583739),
    @(44453, 0, 1, 37.46721618583739),
    @(44454, 0, 1, 37.46721618583739),
    @(44455, 0, 1, 37.46721618583739),
    @(44456, 2, 3, 112.4016485575122),
    @(44457, 0, 3, 112.4016485575122),
    @(44458, 0, 3, 112.4016485575122),
    @(44459, 0, 2, 74.93443237167479)
)

# Copy the date-column formatting (style index used on column A) down from
# the last existing row so the new A-cells carry the same style as the rest
# of the column.
$ws.Cells.Item($lastRow, 1).Copy()

for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = $lastRow + 1 + $i
    $row = $newData[$i]

    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

$excel.CutCopyMode = 0
